$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared strings must be created in the same order as the target
# workbook so the shared string table indices line up with the diff
# (Gulp note first, then the MAINT bug-fix notes). ---
$ws.Range("B28").Value = "Learn Gulp tool"
$ws.Range("B24").Value = "Learn fixing bug: MAINT-3202"
$ws.Range("B25").Value = "Learn fixing bug: MAINT-3222"
$ws.Range("B26").Value = "Learn fixing bug: MAINT-3222"
$ws.Range("B27").Value = "Learn fixing bug: MAINT-3222 and MAINT-3205"

# --- Copy the existing date-formatted cell's formatting down onto the five
# new day cells so they get the same date number format/style as the rest
# of the "Day" column. ---
$ws.Range("A23").Copy()
$ws.Range("A24:A28").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the new week of activity-log rows (Mon 2015-06-15 through
# Fri 2015-06-19). ---
$ws.Range("A24").Value = 42170
$ws.Range("A25").Value = 42171
$ws.Range("A26").Value = 42172
$ws.Range("A27").Value = 42173
$ws.Range("A28").Value = 42174

# --- Scroll the sheet down so the newly added rows are in view and select
# the last entered cell, matching the author's final view/selection state. ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B28").Select()
